# Apply the commit: rename header labels on the existing sheets and add a
# new "PO Forecast" sheet with forecast data (ds / PO_Forecast / yhat_lower /
# yhat_upper).
#
# NOTE: worksheet references handed out by this host are index-based, not
# stable object handles. Worksheets.Add() always inserts the new sheet at
# the front (shifting every other sheet's index by one) and Move() likewise
# reshuffles indices — so ANY worksheet variable captured before such a
# structural change becomes stale (silently points at the wrong sheet)
# afterwards. To stay safe we do all structural changes (add + reorder)
# first, and only then fetch-by-name immediately before each read/write.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "PO Forecast" sheet, then move it after "Monthly Trend"
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsForecast.Move($null, $wsMonthly)

# --- 2. Rename header cells on the existing sheets --------------------------
# Re-fetch fresh references now that no more structural changes will happen.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Populate the new "PO Forecast" sheet --------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Match the page margins used on the other sheets (new sheets default to a
# narrower 0.7/0.7/0.75/0.75/0.3/0.3in layout).
$wsForecast.PageSetup.LeftMargin = 54    # 0.75in
$wsForecast.PageSetup.RightMargin = 54   # 0.75in
$wsForecast.PageSetup.TopMargin = 72     # 1in
$wsForecast.PageSetup.BottomMargin = 72  # 1in
$wsForecast.PageSetup.HeaderMargin = 36  # 0.5in
$wsForecast.PageSetup.FooterMargin = 36  # 0.5in

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold / centered / bordered header look used on the other sheets.
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous (thin box around each cell)

$data = @(
    @(45333.99999999999, 115, 29.07704049264214,  191.2258935531436),
    @(45347.99999999999, 103, 26.54348589433834,  185.2084335605248),
    @(45361.99999999999,  92, 7.512602937479897,  171.1152464075013),
    @(45382.99999999999,  75, -7.060619302550235, 161.5270810883714),
    @(45396.99999999999,  64, -18.2004415610809,  145.7932776904268),
    @(45410.99999999999,  52, -27.21709778338595, 138.121270757534),
    @(45431.99999999999,  35, -39.61854732264456, 119.5852762229287),
    @(45515.99999999999,   0, -112.6792077705216, 51.23129680829454),
    @(45522.99999999999,   0, -117.9306563837489, 39.02832308781102),
    @(45529.99999999999,   0, -124.4445558340624, 33.30561507880076),
    @(45536.99999999999,   0, -131.3318866590995, 34.65159813832357),
    @(45543.99999999999,   0, -136.0414485490621, 32.79815682214603),
    @(45550.99999999999,   0, -137.517972666,     22.5527982567921),
    @(45557.99999999999,   0, -148.2199534304436, 13.22711703271167),
    @(45564.99999999999,   0, -153.3316571229016, 10.19563334744183),
    @(45571.99999999999,   0, -163.6497011964348, 1.085709060724422)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match the date-style formatting used on column A of the other sheets.
$wsForecast.Range("A2:A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
